$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.202.40'
$ws.Range("E2").Value = '  +0.21%  '

$ws.Range("D3").Value = '1.858.59'
$ws.Range("E3").Value = '  +0.06%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.09'
$ws.Range("E5").Value = '  +0.90%  '

$ws.Range("E6").Value = '  -0.05%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4670'
$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2852'
$ws.Range("E8").Value = '  +1.22%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06516'
$ws.Range("E9").Value = '  -0.75%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.68'
$ws.Range("E10").Value = '  +8.55%  '

$ws.Range("E11").Value = '  +0.84%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '97.30'
$ws.Range("E12").Value = '  +0.63%  '

$ws.Range("D13").Value = '1.864.52'
$ws.Range("E13").Value = '  +0.34%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.150'
$ws.Range("E14").Value = '  +0.79%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6785'
$ws.Range("E15").Value = '  +1.96%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '278.88'
$ws.Range("E16").Value = '  -1.27%  '

$ws.Range("D17").Value = '30.202.39'
$ws.Range("E17").Value = '  +0.08%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.60'
$ws.Range("E18").Value = '  +7.92%  '

$ws.Range("E19").Value = '  -0.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.368'
$ws.Range("E20").Value = '  -1.22%  '

$ws.Range("D21").Value = '2.109.31'
$ws.Range("E21").Value = '  -0.02%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.000007305'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  -0.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.150'
$ws.Range("E24").Value = '  +0.19%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '166.75'
$ws.Range("E25").Value = '  -0.66%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.183'
$ws.Range("E26").Value = '  -1.51%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.04'
$ws.Range("E27").Value = '  +0.80%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.926'
$ws.Range("E28").Value = '  +0.57%  '

$ws.Range("E29").Value = '  +3.55%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09707'
$ws.Range("E30").Value = '  +1.44%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.361'

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.476'
$ws.Range("E32").Value = '  +0.37%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.035'
$ws.Range("E33").Value = '  -1.64%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04717'
$ws.Range("E34").Value = '  +1.05%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.131'
$ws.Range("E35").Value = '  +2.97%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7048'

$ws.Range("E37").Value = '  +0.30%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01860'
$ws.Range("E38").Value = '  +0.47%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.629'
$ws.Range("E39").Value = '  +4.78%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.335'
$ws.Range("E40").Value = '  +0.08%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.40'
$ws.Range("E41").Value = '  +3.33%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.948'
$ws.Range("E42").Value = '  +1.15%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8464'
$ws.Range("E43").Value = '  -0.80%  '

$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4168'
$ws.Range("E44").Value = '  +0.22%  '

$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9997'
$ws.Range("E45").Value = '  -0.08%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.31'
$ws.Range("E46").Value = '  -0.34%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '977.21'
$ws.Range("E47").Value = '  -1.72%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.164'
$ws.Range("E48").Value = '  -0.97%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.230'
$ws.Range("E49").Value = '  +0.60%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.04'
$ws.Range("E50").Value = '  +0.37%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05636'
$ws.Range("E51").Value = '  +0.13%  '
